$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume-change (E) figures to match the latest
# coinranking.com snapshot. A handful of rows (Maker/OKB, InjectiveProtocol/
# RenderToken, Stellar/Cosmos) swapped rank order, so B/C/D/E are all rewritten
# for those rows.

# Row 2
$ws.Range("D2").Value = "63.540.08"
$ws.Range("E2").Value = "  -1.90%  "

# Row 3
$ws.Range("D3").Value = "3.383.69"
$ws.Range("E3").Value = "  -1.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "3.382.21"
$ws.Range("E8").Value = "  -1.68%  "

# Row 9
$ws.Range("E9").Value = "  -4.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.119"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.421"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.60%  "

# Row 13
$ws.Range("D13").Value = "3.965.80"
$ws.Range("E13").Value = "  -1.58%  "

# Row 14
$ws.Range("E14").Value = "  +0.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.00%  "

# Row 16
$ws.Range("E16").Value = "  -2.52%  "

# Row 17
$ws.Range("D17").Value = "63.581.32"
$ws.Range("E17").Value = "  -1.83%  "

# Row 18
$ws.Range("D18").Value = "3.402.48"
$ws.Range("E18").Value = "  -2.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.62%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.91%  "

# Row 28
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.47%  "

# Row 31
$ws.Range("E31").Value = "  -5.80%  "

# Row 32
$ws.Range("E32").Value = "  -0.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.48%  "

# Row 35
$ws.Range("E35").Value = "  -4.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.853"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0719"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.88%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.70%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.739.02"
$ws.Range("E41").Value = "  -5.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.68%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.53%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "325.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.12%  "

# Row 49
$ws.Range("E49").Value = "  -4.65%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.36%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.102"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.58%  "
